# [Jimmy] merge de cambios de diapositivas de puntos clientes
#
# 1) Rename sheet "Empleados" -> "Clientes"
# 2) Drop the (now redundant) explicit default style on column G
#    - the column carried its own cellXf (a no-op "General"-format xf)
#      that duplicated what the individual cells already specify, so we
#      clear the column-level formatting and restore each cell's own
#      look (which is untouched / taken from sibling cells that already
#      carry the exact same formatting) so the sheet renders identically
#      but the column no longer pins a dedicated style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename sheet ---------------------------------------------------
$ws.Name = "Clientes"

# --- 2) Column G formatting cleanup ------------------------------------
# Remove the column-wide default style (this also resets every cell in
# the column to the workbook default, so we immediately restore the
# original look of each populated cell below).
$col = $ws.Columns.Item(7)
$col.ClearFormats()

# Header cell G1 matches the same header formatting already used by the
# rest of row 1 (bold, centered, wrapped, text format) - copy it from A1
# which already carries that exact style.
$headerSrc = $ws.Range("A1")
$headerSrc.Copy()
$ws.Range("G1").PasteSpecial(-4122)

# Data cells G2:G5 match the same plain text-formatted style already
# used elsewhere in the data rows - copy it from A2 which already
# carries that exact style.
$dataSrc = $ws.Range("A2")
$dataSrc.Copy()
$ws.Range("G2:G5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
